# Update "Programs ranking based on no. of beneficiaries & expenditure.xlsx"
#
# Changes applied:
#  1. Remove the "China / Unemployment social assistance / 2,300,000" row from the
#     ben_actual sheet (it sits inside the Excel table there), which shifts the
#     following row up and shrinks the table range accordingly.
#  2. Make "exp_actual" the active/selected worksheet (was "ben_planned").
#  3. Update the remembered cell selection on a couple of sheets to reflect where
#     the user had last clicked.

$wb = $excel.ActiveWorkbook

# 1. Delete the obsolete "China" row (row 5) from the ben_actual sheet/table.
$wsBenActual = $wb.Worksheets.Item("ben_actual")
[void]$wsBenActual.Rows.Item(5).Delete()

# Leave the selection on ben_actual where the user left it.
[void]$wsBenActual.Activate()
[void]$wsBenActual.Range("B11").Select()

# Leave the selection on exp_planned where the user left it.
$wsExpPlanned = $wb.Worksheets.Item("exp_planned")
[void]$wsExpPlanned.Activate()
[void]$wsExpPlanned.Range("B22").Select()

# 2. Make exp_actual the final active sheet (this also updates workbook.xml's
#    activeTab and moves the tabSelected flag in the sheetViews).
$wsExpActual = $wb.Worksheets.Item("exp_actual")
[void]$wsExpActual.Activate()
